$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A7 value (was "Киров Пугачева 6") and remove rows 8-10
# (which held "Киров Пугачева 7/8/9"), shrinking the used range to A1:A7
$ws.Range("A7").Value = "Киров Пугачева 11111"
$ws.Rows("8:10").Delete()

# Update selection to match target (J13)
$ws.Range("J13").Select() | Out-Null
